$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph (bold label + text run)
#    that currently sits right after the H1 title paragraph.
# ------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Meta description:*") {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2) At the end of the document, insert a new bold paragraph
#    "Play Double Diamond for Free - Classic Online Slot" right
#    before the final (italic "Prompt: DALLE...") paragraph, and
#    replace that final paragraph's text with the new review blurb
#    (keeping its italic run formatting) - done together via a
#    single InsertXML call so the merge lands cleanly.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$insertPoint = $lastPara.Range
$insertPoint.Collapse(1)

$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Double Diamond for Free - Classic Online Slot</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review and play Double Diamond for free at top online casinos. Enjoy the classic three-reel design and Wild symbol payouts.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($xmlFragment)

# InsertXML leaves a spare empty trailing paragraph behind (Word always
# keeps one final paragraph mark) - fold it back into the paragraph
# before it so the document ends exactly on our new italic paragraph.
$newCount = $d.Paragraphs.Count
if ($d.Paragraphs($newCount).Range.Text.Length -le 1) {
    $trailer = $d.Range($d.Paragraphs($newCount - 1).Range.End - 1, $d.Paragraphs($newCount).Range.End)
    $trailer.Delete()
}
